$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 98.57143000000001
$ws.Range("I9").Value = 98.57143000000001
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 98.57143000000001
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 70.42856999999999
$ws.Range("N9").ClearContents()

$ws.Range("H125").Value = 893.6
$ws.Range("I125").Value = 744
$ws.Range("K125").Value = 6696
$ws.Range("M125").Value = -4236

$ws.Range("H132").Value = 1341.4166
$ws.Range("I132").Value = 1376.0333
$ws.Range("K132").Value = 4128.0999
$ws.Range("M132").Value = -1598.0999

$ws.Range("H135").Value = 726.7313
$ws.Range("J135").Value = 1846.4
$ws.Range("L135").Value = 16617.6
$ws.Range("N135").Value = -21687.6

$ws.Range("H137").Value = 787.8222
$ws.Range("I137").Value = 725.6087
$ws.Range("J137").Value = 852.86365
$ws.Range("K137").Value = 2176.8261
$ws.Range("L137").Value = 2558.59095
$ws.Range("M137").Value = 373.1738999999998
$ws.Range("N137").Value = -7658.59095

$ws.Range("H141").Value = 2382.8
$ws.Range("I141").Value = 733.4
$ws.Range("J141").Value = 7331
$ws.Range("K141").Value = 2200.2
$ws.Range("L141").Value = 21993
$ws.Range("M141").Value = 2979.8
$ws.Range("N141").Value = -32353

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1944.01
$ws.Range("I32").Value = 1879.4744
$ws.Range("J32").Value = 2172.818
$ws.Range("K32").Value = 1879.4744
$ws.Range("L32").Value = 2172.818
$ws.Range("M32").Value = -1592.4744
$ws.Range("N32").Value = -2746.818

$ws.Range("H74").Value = 910.7907
$ws.Range("I74").Value = 872.29266
$ws.Range("K74").Value = 872.29266
$ws.Range("M74").Value = 1.707340000000045

$ws.Range("H77").Value = 910.7907
$ws.Range("I77").Value = 872.29266
$ws.Range("K77").Value = 4361.463299999999
$ws.Range("M77").Value = 6.536700000000565

$ws.Range("H132").Value = 1224.3541
$ws.Range("I132").Value = 1130.561
$ws.Range("J132").Value = 1773.7142
$ws.Range("K132").Value = 3391.683
$ws.Range("L132").Value = 5321.142599999999
$ws.Range("M132").Value = -861.683
$ws.Range("N132").Value = -10381.1426

$ws.Range("H139").Value = 80000
$ws.Range("J139").Value = 80000
$ws.Range("L139").Value = 80000
$ws.Range("N139").Value = -90280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 1468.3334
$ws.Range("I10").Value = 1468.3334
$ws.Range("K10").Value = 1468.3334
$ws.Range("M10").Value = -1328.3334

$ws.Range("H134").Value = 19659.482
$ws.Range("I134").Value = 1525.4783
$ws.Range("J134").Value = 103075.9
$ws.Range("K134").Value = 4576.4349
$ws.Range("L134").Value = 309227.7
$ws.Range("M134").Value = -2041.4349
$ws.Range("N134").Value = -314297.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2121.7778
$ws.Range("I31").Value = 1911.3541
$ws.Range("J31").Value = 2542.625
$ws.Range("K31").Value = 1911.3541
$ws.Range("L31").Value = 2542.625
$ws.Range("M31").Value = -1616.3541
$ws.Range("N31").Value = -3132.625

$ws.Range("H34").Value = 2121.7778
$ws.Range("I34").Value = 1911.3541
$ws.Range("J34").Value = 2542.625
$ws.Range("K34").Value = 1911.3541
$ws.Range("L34").Value = 2542.625
$ws.Range("M34").Value = -1709.3541
$ws.Range("N34").Value = -2946.625

$ws.Range("H58").Value = 946.5
$ws.Range("I58").Value = 1083.0714
$ws.Range("J58").Value = 707.5
$ws.Range("K58").Value = 1083.0714
$ws.Range("L58").Value = 707.5
$ws.Range("M58").Value = -880.0714
$ws.Range("N58").Value = -1113.5

$ws.Range("H132").Value = 1461.3135
$ws.Range("I132").Value = 899.5814
$ws.Range("J132").Value = 2467.75
$ws.Range("K132").Value = 2698.7442
$ws.Range("L132").Value = 7403.25
$ws.Range("M132").Value = -168.7442000000001
$ws.Range("N132").Value = -12463.25

$ws.Range("H134").Value = 1124.2054
$ws.Range("I134").Value = 1058.5172
$ws.Range("J134").Value = 1378.2
$ws.Range("K134").Value = 3175.5516
$ws.Range("L134").Value = 4134.6
$ws.Range("M134").Value = -640.5515999999998
$ws.Range("N134").Value = -9204.6

$ws.Range("H136").Value = 946.5
$ws.Range("I136").Value = 1083.0714
$ws.Range("J136").Value = 707.5
$ws.Range("K136").Value = 3249.2142
$ws.Range("L136").Value = 2122.5
$ws.Range("M136").Value = -699.2142000000003
$ws.Range("N136").Value = -7222.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 389.5
$ws.Range("I17").Value = 199
$ws.Range("K17").Value = 597
$ws.Range("M17").Value = -428

$ws.Range("H31").Value = 5000
$ws.Range("J31").Value = 5000
$ws.Range("L31").Value = 15000
$ws.Range("N31").Value = -15576

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 8114783
$ws.Range("I122").Value = 7983103
$ws.Range("J122").Value = 8334250
$ws.Range("K122").Value = 23949309
$ws.Range("L122").Value = 25002750
$ws.Range("M122").Value = -23946859
$ws.Range("N122").Value = -25007650

$ws.Range("H132").Value = 2059.1396
$ws.Range("I132").Value = 2088.7407
$ws.Range("J132").Value = 2009.1875
$ws.Range("K132").Value = 6266.222099999999
$ws.Range("L132").Value = 6027.5625
$ws.Range("M132").Value = -3736.222099999999
$ws.Range("N132").Value = -11087.5625

$ws.Range("H140").Value = 74750
$ws.Range("J140").Value = 74750
$ws.Range("L140").Value = 74750
$ws.Range("N140").Value = -85110

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3834164.5
$ws.Range("I7").Value = 2787.8572
$ws.Range("K7").Value = 2787.8572
$ws.Range("M7").Value = -2675.8572

$ws.Range("H126").Value = 3834164.5
$ws.Range("I126").Value = 2787.8572
$ws.Range("K126").Value = 8363.571599999999
$ws.Range("M126").Value = -5893.571599999999

$ws.Range("H132").Value = 1229.581
$ws.Range("I132").Value = 1181.8358
$ws.Range("J132").Value = 1686.5714
$ws.Range("K132").Value = 3545.5074
$ws.Range("L132").Value = 5059.7142
$ws.Range("M132").Value = -1015.5074
$ws.Range("N132").Value = -10119.7142

$ws.Range("H136").Value = 1609.9231
$ws.Range("I136").Value = 967.8958
$ws.Range("J136").Value = 3422.7058
$ws.Range("K136").Value = 2903.6874
$ws.Range("L136").Value = 10268.1174
$ws.Range("M136").Value = -353.6873999999998
$ws.Range("N136").Value = -15368.1174

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1386.2858
$ws.Range("I122").Value = 1068
$ws.Range("J122").Value = 1625
$ws.Range("K122").Value = 3204
$ws.Range("L122").Value = 4875
$ws.Range("M122").Value = -754
$ws.Range("N122").Value = -9775

$ws.Range("H132").Value = 474.58334
$ws.Range("I132").Value = 422.29166
$ws.Range("K132").Value = 1266.87498
$ws.Range("M132").Value = 1263.12502

$ws.Range("H136").Value = 772.1489
$ws.Range("I136").Value = 919.7586
$ws.Range("J136").Value = 534.3333
$ws.Range("K136").Value = 2759.2758
$ws.Range("L136").Value = 1602.9999
$ws.Range("M136").Value = -209.2757999999999
$ws.Range("N136").Value = -6702.9999

$ws.Range("H141").Value = 91900
$ws.Range("J141").Value = 91900
$ws.Range("L141").Value = 91900
$ws.Range("N141").Value = -102260
